# Append one new data row (row 99) to "Sheet 1", matching the existing
# table's shape: date (col A, same date/time style as the rows above),
# volume/high/low/open/close (cols B-F, numbers), adj_close (col G,
# text "2"), ticker (col H, text "KK.MI").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 99

# Column A: serial date/time value. Copy the number format (and font)
# from the cell above so it reuses the workbook's existing date style
# instead of minting a new one.
$ws.Cells.Item($row, 1).Value = 45456.2916666667
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats

# Columns B-F: plain numeric values.
$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 2
$ws.Cells.Item($row, 4).Value = 2
$ws.Cells.Item($row, 5).Value = 2
$ws.Cells.Item($row, 6).Value = 2

# Column G: stored as text "2" (not the number 2). The leading apostrophe
# forces text entry; resetting the style back to Normal afterwards drops
# the transient "quote prefix" style so the cell stays on the default
# style like the rest of the column.
$ws.Cells.Item($row, 7).Value = "'2"
$ws.Cells.Item($row, 7).Style = "Normal"

# Column H: ticker text.
$ws.Cells.Item($row, 8).Value = "KK.MI"
